# Added code for Create New Alert functionality
# - Adds a new "Alerts" worksheet (after "Tasks")
# - Populates header row + two sample data rows
# - Updates selection/active-tab bookkeeping to match the authored edit

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Create the new "Alerts" worksheet and move it to the end (after Tasks)
# ------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "Alerts"
$newSheet.Move($null, $wb.Worksheets.Item("Tasks"))

# Re-fetch the sheet by name: the reference above can go stale after Move()
$alerts = $wb.Worksheets.Item("Alerts")

# ------------------------------------------------------------------
# 2. Header row
# ------------------------------------------------------------------
$alerts.Range("A1").Value = "AlertSendMode"
$alerts.Range("B1").Value = "Type"
$alerts.Range("C1").Value = "Trigger"
$alerts.Range("D1").Value = "Format"
$alerts.Range("E1").Value = "Active"

$headerRange = $alerts.Range("A1:E1")
$headerRange.Interior.Color = 65535

# ------------------------------------------------------------------
# 3. Data rows - filled in column by column (A down, then B down, ...)
# ------------------------------------------------------------------
$alerts.Range("A2").Value = "[O]wner"
$alerts.Range("A3").Value = "[S]pecific User (select below)"

$alerts.Range("B2").Value = "Contacts"
$alerts.Range("B3").Value = "Events"

$alerts.Range("C2").Value = "Ownership Changed"
$alerts.Range("C3").Value = "Note Added"

$alerts.Range("D2").Value = "Email Alert"
$alerts.Range("D3").Value = "Text Message Alert"

$alerts.Range("E2").Value = "Yes"
$alerts.Range("E3").Value = "No"

# ------------------------------------------------------------------
# 4. Column widths (best-fit-ish)
# ------------------------------------------------------------------
$alerts.Columns.Item(1).ColumnWidth = 24
$alerts.Columns.Item(3).ColumnWidth = 16.5

# ------------------------------------------------------------------
# 5. Tasks sheet: selection moves from I19 to F32, tabSelected is dropped
# ------------------------------------------------------------------
$tasks = $wb.Worksheets.Item("Tasks")
$tasks.Range("F32").Select()

# ------------------------------------------------------------------
# 6. FeedbackForms sheet: header row (row 1) becomes fully selected
# ------------------------------------------------------------------
$feedback = $wb.Worksheets.Item("FeedbackForms")
$feedback.Rows.Item(1).Select()

# ------------------------------------------------------------------
# 7. Alerts becomes the active sheet / tab, with E3 selected
# ------------------------------------------------------------------
$alerts.Range("E3").Select()
$alerts.Activate()
